# Reorders the player roster rows (A2:C19) on the active worksheet.
# The header row (row 1) and the last three rows (Paolo Banchero, Chet
# Holmgren, Jakob Poeltl) stay put; the remaining players are moved so
# that the block starting at Russell Westbrook comes first.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("Russell Westbrook", "PG", "Denver Nuggets"),
    @("Jaylen Brown", "SG,SF", "Boston Celtics"),
    @("Dejounte Murray", "PG,SG", "New Orleans Pelicans"),
    @("Nikola Jokic", "C", "Denver Nuggets"),
    @("Rudy Gobert", "C", "Minnesota Timberwolves"),
    @("Clint Capela", "C", "Atlanta Hawks"),
    @("Jerami Grant", "SF,PF", "Portland Trail Blazers"),
    @("Pascal Siakam", "SF,PF", "Indiana Pacers"),
    @("Deni Avdija", "SF,PF", "Portland Trail Blazers"),
    @("Julian Champagnie", "SF,PF", "San Antonio Spurs"),
    @("Jalen Suggs", "PG,SG", "Orlando Magic"),
    @("Chris Paul", "PG", "San Antonio Spurs"),
    @("Bogdan Bogdanovic", "SG,SF", "Atlanta Hawks"),
    @("Ayo Dosunmu", "SG,SF", "Chicago Bulls"),
    @("Jalen Green", "PG,SG", "Houston Rockets"),
    @("Paolo Banchero", "SF,PF", "Orlando Magic"),
    @("Chet Holmgren", "PF,C", "Oklahoma City Thunder"),
    @("Jakob Poeltl", "C", "Toronto Raptors")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}
